# Updated cryptos list on Wed May 31 14:20:06 UTC 2023 with GitHub Actions
#
# Price (column D) and Volume(1h) (column E) values are stored as plain text
# in this sheet. Values that look like a single decimal number would be
# auto-converted to a numeric type by Excel's normal text-entry parsing, so
# those are written with a leading apostrophe to force them to stay text
# (matching the original inline-string content/type), then ClearFormats()
# removes the "quote prefix" marker style Excel applies for the apostrophe,
# restoring the cell to its original unstyled state. Values that already
# can't parse as a single number (e.g. "27.088.88", which uses '.' as a
# thousands separator) are written as-is, with no style side effects.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.088.88"
$ws.Range("E2").Value = "  -2.71%  "

$ws.Range("D3").Value = "1.865.98"
$ws.Range("E3").Value = "  -2.27%  "

$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").Value = "'306.72"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.97%  "

$ws.Range("E6").Value = "  +0.04%  "

$ws.Range("D7").Value = "'0.5122"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +2.37%  "

$ws.Range("D8").Value = "'0.3742"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.22%  "

$ws.Range("D9").Value = "'0.07118"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.18%  "

$ws.Range("D10").Value = "'0.8864"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.92%  "

$ws.Range("E11").Value = "  -3.05%  "

$ws.Range("D12").Value = "1.870.14"
$ws.Range("E12").Value = "  -2.10%  "

$ws.Range("D13").Value = "'0.07541"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.22%  "

$ws.Range("D14").Value = "'5.319"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.79%  "

$ws.Range("D15").Value = "'88.84"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -3.71%  "

$ws.Range("E16").Value = "  -0.15%  "

$ws.Range("D17").Value = "'0.000008462"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -2.95%  "

$ws.Range("E18").Value = "  -3.33%  "

$ws.Range("D19").Value = "'1.000"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.10%  "

$ws.Range("D20").Value = "27.139.74"
$ws.Range("E20").Value = "  -2.62%  "

$ws.Range("D21").Value = "'5.047"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.34%  "

$ws.Range("D22").Value = "2.095.40"
$ws.Range("E22").Value = "  -2.61%  "

$ws.Range("D23").Value = "'10.53"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.90%  "

$ws.Range("D24").Value = "'6.454"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.29%  "

$ws.Range("D25").Value = "'149.56"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.17%  "

$ws.Range("D26").Value = "'1.842"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.08%  "

$ws.Range("D27").Value = "'17.95"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.31%  "

$ws.Range("D28").Value = "'2.086"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -6.25%  "

$ws.Range("D29").Value = "'112.99"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.81%  "

$ws.Range("D30").Value = "'4.695"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -3.93%  "

$ws.Range("D31").Value = "'4.659"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.18%  "

$ws.Range("D32").Value = "'0.09036"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.66%  "

$ws.Range("D33").Value = "'0.05133"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.92%  "

$ws.Range("D34").Value = "'3.076"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.15%  "

$ws.Range("D35").Value = "'1.154"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -6.50%  "

$ws.Range("D36").Value = "'0.7328"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -6.93%  "

$ws.Range("D37").Value = "'0.02055"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.17%  "

$ws.Range("D38").Value = "'2.500"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -5.80%  "

$ws.Range("D39").Value = "'3.056"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.26%  "

$ws.Range("E40").Value = "  -1.99%  "

$ws.Range("D41").Value = "'0.5323"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -3.47%  "

$ws.Range("D42").Value = "'6.579"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.98%  "

$ws.Range("D43").Value = "'116.08"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.49%  "

$ws.Range("D44").Value = "'8.308"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.07%  "

$ws.Range("D45").Value = "'0.1467"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.96%  "

$ws.Range("D46").Value = "'1.001"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.09%  "

$ws.Range("D47").Value = "'0.4611"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.80%  "

$ws.Range("D48").Value = "'10.04"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -5.40%  "

$ws.Range("D49").Value = "'1.564"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -4.32%  "

# Row 50 and 51 swap: Elrond/Aave trade places with updated price/volume.
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "'64.24"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -4.46%  "

$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").Value = "'36.60"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.92%  "
